# Update cached market-price / profit figures on the per-job "Profits" sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect a refreshed price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 722.7143
$ws.Range("I11").Value = 722.7143
$ws.Range("K11").Value = 722.7143
$ws.Range("M11").Value = -582.7143

# Row 62
$ws.Range("H62").Value = 62518516
$ws.Range("I62").Value = 100002136
$ws.Range("J62").Value = 45820.668
$ws.Range("K62").Value = 100002136
$ws.Range("L62").Value = 45820.668
$ws.Range("M62").Value = -100001512
$ws.Range("N62").Value = -47068.668

# Row 65
$ws.Range("H65").Value = 62518516
$ws.Range("I65").Value = 100002136
$ws.Range("J65").Value = 45820.668
$ws.Range("K65").Value = 500010680
$ws.Range("L65").Value = 229103.34
$ws.Range("M65").Value = -500007560
$ws.Range("N65").Value = -235343.34

# Row 98
$ws.Range("H98").Value = 9966.223
$ws.Range("I98").Value = 10025.625
$ws.Range("J98").Value = 9491
$ws.Range("K98").Value = 10025.625
$ws.Range("L98").Value = 9491
$ws.Range("M98").Value = -8527.625
$ws.Range("N98").Value = -12487

# Row 116
$ws.Range("H116").Value = 35724350
$ws.Range("I116").Value = 62506372
$ws.Range("J116").Value = 14996.667
$ws.Range("K116").Value = 62506372
$ws.Range("L116").Value = 14996.667
$ws.Range("M116").Value = -62502930
$ws.Range("N116").Value = -21880.667

# Row 122
$ws.Range("H122").Value = 9966.223
$ws.Range("I122").Value = 10025.625
$ws.Range("J122").Value = 9491
$ws.Range("K122").Value = 30076.875
$ws.Range("L122").Value = 28473
$ws.Range("M122").Value = -27626.875
$ws.Range("N122").Value = -33373

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 3280.6296
$ws.Range("I2").Value = 2136.5625
$ws.Range("J2").Value = 4944.727
$ws.Range("K2").Value = 2136.5625
$ws.Range("L2").Value = 4944.727
$ws.Range("M2").Value = -2023.5625
$ws.Range("N2").Value = -5170.727

# Row 45
$ws.Range("H45").Value = 6217.857
$ws.Range("I45").Value = 2100
$ws.Range("J45").Value = 7865
$ws.Range("K45").Value = 2100
$ws.Range("L45").Value = 7865
$ws.Range("M45").Value = -1723
$ws.Range("N45").Value = -8619

# Row 61
$ws.Range("H61").Value = 7235.2896
$ws.Range("I61").Value = 3403.68
$ws.Range("K61").Value = 3403.68
$ws.Range("M61").Value = -3191.68

# Row 74
$ws.Range("H74").Value = 45825.46
$ws.Range("I74").Value = 58896.82
$ws.Range("K74").Value = 58896.82
$ws.Range("M74").Value = -58022.82

# Row 77
$ws.Range("H77").Value = 45825.46
$ws.Range("I77").Value = 58896.82
$ws.Range("K77").Value = 294484.1
$ws.Range("M77").Value = -290116.1

# Row 102
$ws.Range("H102").Value = 975.3913
$ws.Range("I102").Value = 979
$ws.Range("K102").Value = 979
$ws.Range("M102").Value = 643

# Row 116
$ws.Range("H116").Value = 3280.6296
$ws.Range("I116").Value = 2136.5625
$ws.Range("J116").Value = 4944.727
$ws.Range("K116").Value = 2136.5625
$ws.Range("L116").Value = 4944.727
$ws.Range("M116").Value = 157.4375
$ws.Range("N116").Value = -9532.726999999999

# Row 132
$ws.Range("H132").Value = 3693.5625
$ws.Range("I132").Value = 1480.7059
$ws.Range("K132").Value = 4442.1177
$ws.Range("M132").Value = -1912.1177

# Row 136
$ws.Range("H136").Value = 7235.2896
$ws.Range("I136").Value = 3403.68
$ws.Range("K136").Value = 10211.04
$ws.Range("M136").Value = -7661.039999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 3280.6296
$ws.Range("I3").Value = 2136.5625
$ws.Range("J3").Value = 4944.727
$ws.Range("K3").Value = 2136.5625
$ws.Range("L3").Value = 4944.727
$ws.Range("M3").Value = -2022.5625
$ws.Range("N3").Value = -5172.727

# Row 81
$ws.Range("H81").Value = 78568.28999999999
$ws.Range("J81").Value = 78568.28999999999
$ws.Range("L81").Value = 78568.28999999999
$ws.Range("N81").Value = -80690.28999999999

# Row 84
$ws.Range("H84").Value = 78568.28999999999
$ws.Range("J84").Value = 78568.28999999999
$ws.Range("L84").Value = 235704.87
$ws.Range("N84").Value = -246312.87

# Row 105
$ws.Range("H105").Value = 3057.8
$ws.Range("J105").Value = 3533.2222
$ws.Range("L105").Value = 3533.2222
$ws.Range("N105").Value = -7027.2222

# Row 107
$ws.Range("H107").Value = 41669744
$ws.Range("I107").Value = 48915076
$ws.Range("J107").Value = 9075.25
$ws.Range("K107").Value = 48915076
$ws.Range("L107").Value = 9075.25
$ws.Range("M107").Value = -48913156
$ws.Range("N107").Value = -12915.25

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 6716.294
$ws.Range("I16").Value = 4457.2856
$ws.Range("K16").Value = 4457.2856
$ws.Range("M16").Value = -4170.2856

# Row 31
$ws.Range("H31").Value = 11232.75
$ws.Range("I31").Value = 5475.364
$ws.Range("J31").Value = 14248.523
$ws.Range("K31").Value = 5475.364
$ws.Range("L31").Value = 14248.523
$ws.Range("M31").Value = -5180.364
$ws.Range("N31").Value = -14838.523

# Row 34
$ws.Range("H34").Value = 11232.75
$ws.Range("I34").Value = 5475.364
$ws.Range("J34").Value = 14248.523
$ws.Range("K34").Value = 5475.364
$ws.Range("L34").Value = 14248.523
$ws.Range("M34").Value = -5273.364
$ws.Range("N34").Value = -14652.523

# Row 113
$ws.Range("H113").Value = 6716.294
$ws.Range("I113").Value = 4457.2856
$ws.Range("K113").Value = 4457.2856
$ws.Range("M113").Value = -2287.2856

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 5718844
$ws.Range("I5").Value = 13334003
$ws.Range("K5").Value = 40002009
$ws.Range("M5").Value = -40001897

# Row 39
$ws.Range("H39").Value = 12581.833
$ws.Range("J39").Value = 12581.833
$ws.Range("L39").Value = 37745.499
$ws.Range("N39").Value = -38333.499

# Row 131
$ws.Range("H131").Value = 1759.2858
$ws.Range("J131").Value = 3103.5
$ws.Range("L131").Value = 9310.5
$ws.Range("N131").Value = -19390.5

# Row 135
$ws.Range("H135").Value = 5718844
$ws.Range("I135").Value = 13334003
$ws.Range("K135").Value = 120006027
$ws.Range("M135").Value = -120003492

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 7778.5713
$ws.Range("I113").Value = 3928.5715
$ws.Range("K113").Value = 3928.5715
$ws.Range("M113").Value = -1758.5715

# Row 122
$ws.Range("H122").Value = 2015820.5
$ws.Range("I122").Value = 3152715.2
$ws.Range("J122").Value = 4391.385
$ws.Range("K122").Value = 9458145.600000001
$ws.Range("L122").Value = 13174.155
$ws.Range("M122").Value = -9455695.600000001
$ws.Range("N122").Value = -18074.155

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 2003.3846
$ws.Range("J46").Value = 2544.4
$ws.Range("L46").Value = 2544.4
$ws.Range("N46").Value = -2920.4

# Row 93
$ws.Range("H93").Value = 4810.125
$ws.Range("I93").Value = 3253.5
$ws.Range("J93").Value = 12593.25
$ws.Range("K93").Value = 3253.5
$ws.Range("L93").Value = 12593.25
$ws.Range("M93").Value = -2005.5
$ws.Range("N93").Value = -15089.25

# Row 100
$ws.Range("H100").Value = 4601.1
$ws.Range("I100").Value = 3251.1667
$ws.Range("K100").Value = 3251.1667
$ws.Range("M100").Value = -2710.1667

# Row 122
$ws.Range("H122").Value = 7409.5454
$ws.Range("J122").Value = 8167.222
$ws.Range("L122").Value = 24501.666
$ws.Range("N122").Value = -29401.666

$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Range("H46").Value = 89664.664
$ws.Range("J46").Value = 89664.664
$ws.Range("L46").Value = 89664.664
$ws.Range("N46").Value = -90126.664

# Row 134
$ws.Range("H134").Value = 89664.664
$ws.Range("J134").Value = 89664.664
$ws.Range("L134").Value = 268993.992
$ws.Range("N134").Value = -274063.992

# Row 136
$ws.Range("H136").Value = 43485900
$ws.Range("J136").Value = 9535.588
$ws.Range("L136").Value = 28606.764
$ws.Range("N136").Value = -33706.764
